$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 21:18"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5777573
$ws.Range("C4").Value = 31301
$ws.Range("D4").Value = 3105570
$ws.Range("E4").Value = 2493933
$ws.Range("G4").Value = 646
$ws.Range("H4").Value = 178070

# Espana (row 12)
$ws.Range("B12").Value = 407879
$ws.Range("C12").Value = 3650
$ws.Range("G12").Value = 25
$ws.Range("H12").Value = 28838

# Alemania (row 23)
$ws.Range("B23").Value = 231830
$ws.Range("C23").Value = 546
$ws.Range("E23").Value = 16703
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 9327

# Rows 44/45: Guatemala and Emiratos Arabes Unidos swap rank order
$ws.Range("A44").Value = "Guatemala"
$ws.Range("B44").Value = 66941
$ws.Range("C44").Value = 958
$ws.Range("D44").Value = 55314
$ws.Range("E44").Value = 9095
$ws.Range("G44").Value = 26
$ws.Range("H44").Value = 2532

$ws.Range("A45").Value = "Emiratos Arabes Unidos"
$ws.Range("B45").Value = 66193
$ws.Range("C45").Value = 391
$ws.Range("D45").Value = 58296
$ws.Range("E45").Value = 7527
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 370

# Ghana (row 55)
$ws.Range("B55").Value = 43325
$ws.Range("C55").Value = 65
$ws.Range("D55").Value = 41408
$ws.Range("E55").Value = 1656

# Rows 60/61: Uzbekistan and Afganistan swap rank order
$ws.Range("A60").Value = "Uzbekistan"
$ws.Range("B60").Value = 38074
$ws.Range("C60").Value = 527
$ws.Range("D60").Value = 33989
$ws.Range("E60").Value = 3825
$ws.Range("G60").Value = 8
$ws.Range("H60").Value = 260

$ws.Range("A61").Value = "Afganistan"
$ws.Range("B61").Value = 37894
$ws.Range("C61").Value = 38
$ws.Range("D61").Value = 28016
$ws.Range("E61").Value = 8493
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1385

# Maldivas (row 104)
$ws.Range("E104").Value = 2526
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 26

# Zimbabue (row 105)
$ws.Range("B105").Value = 5815
$ws.Range("C105").Value = 70
$ws.Range("D105").Value = 4587
$ws.Range("E105").Value = 1076
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 152

# Sudan del Sur (row 131)
$ws.Range("B131").Value = 2497
$ws.Range("C131").Value = 3
$ws.Range("E131").Value = 1160

# Republica del Chad (row 161)
$ws.Range("B161").Value = 981
$ws.Range("C161").Value = 9
$ws.Range("E161").Value = 36

# Burundi (row 171)
$ws.Range("B171").Value = 426
$ws.Range("C171").Value = 4
$ws.Range("E171").Value = 89

# Barbados (row 188)
$ws.Range("B188").Value = 157
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 124
